# Automatische test-sync: 2025-06-19 17:48:21
#
# Appends three new mail-log rows to "Logs", inserts a new "Sollicitatie"
# category row into the "Dashboard" summary table (bumping "Afmelding" from
# 4 to 5), and updates the dashboard bar chart's category/value series
# references to include the newly inserted row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append rows 20-22
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A20").Value = "Sollicitatie marketingfunctie"
$logs.Range("B20").Value = "mailmind.test@zohomail.eu"
$logs.Range("C20").Value = "Hierbij solliciteer ik voor de functie van marketeer. Zie bijlage voor CV."
$logs.Range("D20").Value = "Sollicitatie"
$logs.Range("F20").Value = "2025-06-19 17:47:20"
$logs.Range("G20").Value = "Nee"

$logs.Range("A21").Value = "Sollicitatie marketingfunctie"
$logs.Range("B21").Value = "mailmind.test@zohomail.eu"
$logs.Range("C21").Value = "Hierbij solliciteer ik voor de functie van marketeer. Zie bijlage voor CV."
$logs.Range("D21").Value = "Sollicitatie"
$logs.Range("F21").Value = "2025-06-19 17:47:21"
$logs.Range("G21").Value = "Nee"

$logs.Range("A22").Value = "Afmelding nieuwsbrief"
$logs.Range("B22").Value = "mailmind.test@zohomail.eu"
$logs.Range("C22").Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$logs.Range("D22").Value = "Afmelding"
$logs.Range("F22").Value = "2025-06-19 17:48:20"
$logs.Range("G22").Value = "Nee"

# Conditional formatting ranges need to grow from row 19 to row 22 to keep
# covering the newly-added rows.
$logs.Range("D2:D19").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D22"))
$logs.Range("G2:G19").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G22"))

# ---------------------------------------------------------------------
# 2. Dashboard sheet: insert a "Sollicitatie" row before "Informatieaanvraag"
#    and bump the "Afmelding" count.
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("B3").Value = 5

$dash.Rows.Item(6).Insert()
$dash.Range("A6").Value = "Sollicitatie"
$dash.Range("B6").Value = 2

# ---------------------------------------------------------------------
# 3. Chart: extend the category/value series references to the new
#    Dashboard row (A2:A6 -> A2:A7, B2:B6 -> B2:B7). Keep the series name
#    reference (tx) untouched by only rewriting categories/values, and
#    keep the sheet name quoted so the stored ref form matches the
#    original 'Dashboard'!... style.
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$series = $chartObj.Chart.SeriesCollection(1)
$series.Values = "='Dashboard'!`$B`$2:`$B`$7"
$series.XValues = "='Dashboard'!`$A`$2:`$A`$7"
